$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), copying H1's formatting (bold, border, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I (I0) and J (IF) columns, rows 2-18
$iValues = @(4, 2, 3, 4, 5, 3, 6, 5, 6, 5, 6, 4, 1, 6, 1, 1, 2)
$jValues = @(5, 4, 4, 5, 7, 5, 7, 6, 7, 6, 7, 8, 2, 8, 5, 1, 2)

for ($r = 0; $r -lt 17; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
